$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 31 (header row, style like row 25/28 header -> fill/yellow style "s2") ---
$ws.Range("A31").Value = "TC-118"
$ws.Range("B31").Value = "username"
$ws.Range("C31").Value = "password"
$ws.Range("D31").Value = "location"
$ws.Range("E31").Value = "hotel"
$ws.Range("F31").Value = "Check In Date"
$ws.Range("G31").Value = "Room Type"
$ws.Range("H31").Value = "No of rooms"
$ws.Range("I31").Value = "Check out Date"
$ws.Range("J31").Value = "Adults per Room"
$ws.Range("K31").Value = "Children per Room  "
$ws.Range("L31").Value = "Login page title"
$ws.Range("M31").Value = "Hotel search title"
$ws.Range("N31").Value = "Select page title"
$ws.Range("O31").Value = "Booking page title"
$ws.Range("P31").Value = "Logout page title"
$ws.Range("Q31").Value = "message"
$ws.Range("R31").Value = "First Name"
$ws.Range("S31").Value = "LastName"
$ws.Range("T31").Value = "Address"
$ws.Range("U31").Value = "Credit num"
$ws.Range("V31").Value = "Card type"
$ws.Range("W31").Value = "Expiry month"
$ws.Range("X31").Value = "Expiry year"
$ws.Range("Y31").Value = "CVV number"
$ws.Range("Z31").Value = "Hotel select title"
$ws.Range("AA31").Value = "Results found"

# apply the same formatting as the row 25 header (yellow fill header style) to the new header row
$ws.Range("A25:AA25").Copy()
$ws.Range("A31:AA31").PasteSpecial(-4122)

# --- Row 32 (data row, mirrors row 26 pattern) ---
$ws.Range("A32").Value = "TC-118"
$ws.Range("B32").Value = "reyaz0806"
$ws.Range("C32").Value = "reyaz123"
$ws.Range("D32").Value = "Sydney"
$ws.Range("E32").Value = "Hotel Creek"
$ws.Range("F32").Value = "19/01/2025"
$ws.Range("G32").Value = "Standard"
$ws.Range("H32").Value = "1 - One"
$ws.Range("I32").Value = "20/01/2025"
$ws.Range("J32").Value = "1 - One"
$ws.Range("K32").Value = "0 - None"
$ws.Range("L32").Value = "Adactin.com - Hotel Reservation System"
$ws.Range("M32").Value = "Adactin.com - Search Hotel"
$ws.Range("N32").Value = "Adactin.com - Select Hotel"
$ws.Range("O32").Value = "Adactin.com - Book A Hotel"
$ws.Range("P32").Value = "Adactin.com - Logout"
$ws.Range("Q32").Value = "Cancel Selected"
$ws.Range("R32").Value = "Test"
$ws.Range("S32").Value = "Data"
$ws.Range("T32").Value = "Hyderabad"
$ws.Range("U32").Value = "1234567812345678"
$ws.Range("V32").Value = "Master Card"
$ws.Range("W32").Value = "March"
$ws.Range("X32").Value = "2026"
$ws.Range("Y32").Value = "000"
$ws.Range("Z32").Value = "Adactin.com - Select Hotel"
$ws.Range("AA32").Value = "1 result(s) found. Show all"

# apply the same per-cell formatting as row 26 (plain cells + date/quote-prefix/wrap styles)
$ws.Range("A26:AA26").Copy()
$ws.Range("A32:AA32").PasteSpecial(-4122)

# row 32 wraps text in column Q (like row 26/23/20) -> taller row
$ws.Rows.Item(32).RowHeight = 28.8

# --- column widths for the newly used columns Y (25) and Z (26) ---
$ws.Columns.Item(25).ColumnWidth = 12.109375
$ws.Columns.Item(26).ColumnWidth = 17.6640625

# --- update selection to match the new active cell ---
$ws.Range("AE32").Select()
